# The underlying distance-matrix data didn't change, only two Turkish city
# labels were shortened (their parenthetical "also-known-as" qualifier was
# dropped):
#   "KOCAELİ (İZMİT)"      -> "KOCAELİ"   (column R header, row 18 label)
#   "SAKARYA (ADAPAZARI)"  -> "SAKARYA"   (column Z header, row 26 label)
#
# Each city name appears twice: once as a column header in row 1 and once as
# the row label in column A (since the sheet is a symmetric distance matrix).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column headers (row 1)
$ws.Range("R1").Value = "KOCAELİ"
$ws.Range("Z1").Value = "SAKARYA"

# Row labels (column A)
$ws.Range("A18").Value = "KOCAELİ"
$ws.Range("A26").Value = "SAKARYA"

# The header row auto-wraps/rotates its text; with shorter labels the row
# needs less vertical space.
$ws.Rows.Item(1).RowHeight = 84.75

# The active selection moved to Z2.
$ws.Activate()
$ws.Range("Z2").Select()
